$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.919.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.88%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.509.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.11%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.14%  '
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.508'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.505.29'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.71%  '
$ws.Range("E11").Value = '  -1.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.342'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.990.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '69.928.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000178'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.520.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.47%  '
$ws.Range("E23").Value = '  -3.70%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '68.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.642.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0905'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '477.51'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.28'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("E34").Value = '  -4.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '156.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.49'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.64%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.29'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.318'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.59'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.521'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.595'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.33%  '
